# Update the document's style sheet (word/styles.xml):
#   - the "east Asian" font used by the Normal and Heading paragraph
#     styles changes from "DejaVu Sans" to "Tahoma"
#   - the List, Caption and Index paragraph styles pick up an explicit
#     complex-script ("cs") font of "DejaVu Sans" on their <w:rPr>
#
# Word's Style object exposes a Font sub-object whose NameFarEast /
# NameBi properties round-trip to <w:rFonts w:eastAsia="…"/> and
# <w:rFonts w:cs="…"/> respectively, so the style definitions can be
# edited directly through Styles(name).Font just like character
# formatting on a Range.

$d = $word.ActiveDocument

# east Asian font: DejaVu Sans -> Tahoma
$d.Styles("Normal").Font.NameFarEast  = "Tahoma"
$d.Styles("Heading").Font.NameFarEast = "Tahoma"

# explicit complex-script font on styles that previously had none
$d.Styles("List").Font.NameBi    = "DejaVu Sans"
$d.Styles("Caption").Font.NameBi = "DejaVu Sans"
$d.Styles("Index").Font.NameBi   = "DejaVu Sans"
